$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 191.89
$ws.Range("I15").Value = 191.89
$ws.Range("K15").Value = 575.67
$ws.Range("M15").Value = -406.67

$ws.Range("H62").Value = 2609
$ws.Range("I62").Value = 2533.2222
$ws.Range("J62").Value = 2950
$ws.Range("K62").Value = 2533.2222
$ws.Range("L62").Value = 2950
$ws.Range("M62").Value = -1909.2222
$ws.Range("N62").Value = -4198

$ws.Range("H64").Value = 69691.53
$ws.Range("I64").Value = 169166.67
$ws.Range("K64").Value = 169166.67
$ws.Range("M64").Value = -168918.67

$ws.Range("H65").Value = 2609
$ws.Range("I65").Value = 2533.2222
$ws.Range("J65").Value = 2950
$ws.Range("K65").Value = 12666.111
$ws.Range("L65").Value = 14750
$ws.Range("M65").Value = -9546.111
$ws.Range("N65").Value = -20990

$ws.Range("H67").Value = 69691.53
$ws.Range("I67").Value = 169166.67
$ws.Range("K67").Value = 169166.67
$ws.Range("M67").Value = -168308.67

$ws.Range("H103").Value = 4307
$ws.Range("J103").Value = 10001
$ws.Range("L103").Value = 30003
$ws.Range("N103").Value = -31175

$ws.Range("H129").Value = 3210.628
$ws.Range("I129").Value = 33785
$ws.Range("J129").Value = 917.55
$ws.Range("K129").Value = 101355
$ws.Range("L129").Value = 2752.65
$ws.Range("M129").Value = -96355
$ws.Range("N129").Value = -12752.65

$ws.Range("H138").Value = 2574.77
$ws.Range("I138").Value = 1185.45
$ws.Range("J138").Value = 2922.1
$ws.Range("K138").Value = 3556.35
$ws.Range("L138").Value = 8766.3
$ws.Range("M138").Value = 1583.65
$ws.Range("N138").Value = -19046.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9841.75
$ws.Range("I32").Value = 8219.658
$ws.Range("J32").Value = 25251.625
$ws.Range("K32").Value = 8219.658
$ws.Range("L32").Value = 25251.625
$ws.Range("M32").Value = -7932.657999999999
$ws.Range("N32").Value = -25825.625

$ws.Range("H74").Value = 1141.2
$ws.Range("I74").Value = 1146.3334
$ws.Range("K74").Value = 1146.3334
$ws.Range("M74").Value = -272.3334

$ws.Range("H77").Value = 1141.2
$ws.Range("I77").Value = 1146.3334
$ws.Range("K77").Value = 5731.666999999999
$ws.Range("M77").Value = -1363.666999999999

$ws.Range("H122").Value = 2222.2856
$ws.Range("I122").Value = 1795.4736
$ws.Range("K122").Value = 5386.4208
$ws.Range("M122").Value = -2936.4208

$ws.Range("H132").Value = 13831.039
$ws.Range("I132").Value = 19409.47
$ws.Range("J132").Value = 2674.1765
$ws.Range("K132").Value = 58228.41
$ws.Range("L132").Value = 8022.529500000001
$ws.Range("M132").Value = -55698.41
$ws.Range("N132").Value = -13082.5295

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 43188.6
$ws.Range("I20").Value = 50981.24
$ws.Range("J20").Value = 2277.25
$ws.Range("K20").Value = 50981.24
$ws.Range("L20").Value = 2277.25
$ws.Range("M20").Value = -50734.24
$ws.Range("N20").Value = -2771.25

$ws.Range("H64").Value = 3000.1667
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 3000.1667
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 3000.1667
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -3450.1667

$ws.Range("H67").Value = 3000.1667
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 3000.1667
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 3000.1667
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -4560.1667

$ws.Range("H107").Value = 47641972
$ws.Range("I107").Value = 58851500
$ws.Range("K107").Value = 58851500
$ws.Range("M107").Value = -58849580

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 25167.475
$ws.Range("I31").Value = 791.08105
$ws.Range("J31").Value = 62747.75
$ws.Range("K31").Value = 791.08105
$ws.Range("L31").Value = 62747.75
$ws.Range("M31").Value = -496.08105
$ws.Range("N31").Value = -63337.75

$ws.Range("H34").Value = 25167.475
$ws.Range("I34").Value = 791.08105
$ws.Range("J34").Value = 62747.75
$ws.Range("K34").Value = 791.08105
$ws.Range("L34").Value = 62747.75
$ws.Range("M34").Value = -589.08105
$ws.Range("N34").Value = -63151.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 75
$ws.Range("J60").Value = 100
$ws.Range("L60").Value = 300
$ws.Range("N60").Value = -802

$ws.Range("H87").Value = 13000
$ws.Range("I87").Value = 13000
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 39000
$ws.Range("L87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value = -37752

$ws.Range("H90").Value = 13000
$ws.Range("I90").Value = 13000
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 117000
$ws.Range("L90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("N90").Value = -110760

$ws.Range("H98").Value = 92005.27
$ws.Range("J98").Value = 92005.27
$ws.Range("L98").Value = 276015.81
$ws.Range("N98").Value = -279011.81

$ws.Range("H107").Value = 1137.7858
$ws.Range("I107").Value = 528.7059
$ws.Range("J107").Value = 2079.0908
$ws.Range("K107").Value = 1586.1177
$ws.Range("L107").Value = 6237.2724
$ws.Range("M107").Value = 333.8822999999998
$ws.Range("N107").Value = -10077.2724

$ws.Range("H113").Value = 794.5405
$ws.Range("J113").Value = 533.2
$ws.Range("L113").Value = 1599.6
$ws.Range("N113").Value = -5939.6

$ws.Range("H132").Value = 418421.78
$ws.Range("I132").Value = 810.36365
$ws.Range("J132").Value = 771785.3
$ws.Range("K132").Value = 7293.27285
$ws.Range("L132").Value = 6946067.7
$ws.Range("M132").Value = -4763.27285
$ws.Range("N132").Value = -6951127.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 216383.39
$ws.Range("I102").Value = 1986.75
$ws.Range("J102").Value = 502245.6
$ws.Range("K102").Value = 1986.75
$ws.Range("L102").Value = 502245.6
$ws.Range("M102").Value = -364.75
$ws.Range("N102").Value = -505489.6

$ws.Range("H122").Value = 3826
$ws.Range("I122").Value = 2933.3333
$ws.Range("K122").Value = 8799.999899999999
$ws.Range("M122").Value = -6349.999899999999

$ws.Range("H132").Value = 2649.8125
$ws.Range("I132").Value = 1739
$ws.Range("J132").Value = 4977.4443
$ws.Range("K132").Value = 5217
$ws.Range("L132").Value = 14932.3329
$ws.Range("M132").Value = -2687
$ws.Range("N132").Value = -19992.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 84329.164
$ws.Range("I100").Value = 167344.67
$ws.Range("J100").Value = 1313.6666
$ws.Range("K100").Value = 334689.34
$ws.Range("L100").Value = 2627.3332
$ws.Range("M100").Value = -334148.34
$ws.Range("N100").Value = -3709.3332

$ws.Range("H107").Value = 111850.336
$ws.Range("J107").Value = 250599.75
$ws.Range("L107").Value = 751799.25
$ws.Range("N107").Value = -755639.25

$ws.Range("H127").Value = 29900
$ws.Range("J127").Value = 29900
$ws.Range("L127").Value = 29900
$ws.Range("N127").Value = -39820

$ws.Range("H132").Value = 3654.59
$ws.Range("I132").Value = 2027.675
$ws.Range("J132").Value = 6753.476
$ws.Range("K132").Value = 6083.025
$ws.Range("L132").Value = 20260.428
$ws.Range("M132").Value = -3553.025
$ws.Range("N132").Value = -25320.428
